$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")
$tr = $shp.TextFrame.TextRange

# Setting the text directly from the original (mis-spelled, multi-run) value
# to the corrected value causes the host to diff run-by-run and keep the
# text split across several runs (preserving the old per-run formatting,
# including the err="1" spell-check flag on the old "PenManship" run).
# Assigning a throwaway placeholder first collapses the paragraph down to a
# single run, so the final assignment lands as one clean run that matches
# the sibling runs' formatting (lang="en-US" dirty="0"), exactly like the
# corrected slide.
$tr.Text = "x"
$tr.Text = "Learn2Write is an interactive learning tool designed to evaluate Penmanship skills using a trained computer model. "
